# Terminal Hortofrutícola Agro Chillán - Acelga
# Weekly update: insert a new latest-week record at row 123 (pushing the
# existing historical rows 123:229 down to 124:230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 123, shifting rows
# 123:229 down to 124:230 (the former row 229 becomes row 230).
$ws.Rows("123:123").Insert()

# Populate the new row 123 with this week's record.
$ws.Range("A123").Value = 7
$ws.Range("B123").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C123").Value = "Ñuble"
$ws.Range("D123").Value = 44669
$ws.Range("E123").Value = 16
$ws.Range("F123").Value = 100112009
$ws.Range("G123").Value = "Acelga"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 200
$ws.Range("K123").Value = 550
$ws.Range("L123").Value = 600
$ws.Range("M123").Value = 575
$ws.Range("N123").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O123").Value = "Provincia de Diguillín"
$ws.Range("P123").Value = 575
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = "Hortaliza"
